{"js": "// Loop 45 (PDF Exporter) \u2014 append the OpenPDF Maven dependency block\n// after the existing dependencies, mirroring the Loop 43/44 snippets\n// already present in the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The target is the last (trailing empty) paragraph in the document body,\n// immediately before the section properties.\nlet anchor = paragraphs.items[paragraphs.items.length - 1];\n\n// Add the Maven-repository comment line into that existing empty paragraph.\nanchor.insertText(\n  \"<!-- https://mvnrepository.com/artifact/com.github.librepdf/openpdf -->\",\n  Word.InsertLocation.end\n);\n\n// Each subsequent line becomes its own new paragraph inserted after the\n// previous one, matching the diff's five inserted <w:p> blocks.\nconst lines = [\n  \"<dependency>\",\n  \"    <groupId>com.github.librepdf</groupId>\",\n  \"    <artifactId>openpdf</artifactId>\",\n  \"    <version>1.3.25</version>\",\n  \"</dependency>\"\n];\n\nfor (const line of lines) {\n  anchor = anchor.insertParagraph(line, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Loop 45 (PDF Exporter) \u2014 append the OpenPDF Maven dependency block\n# after the existing dependencies, mirroring the Loop 43/44 snippets\n# already present in the document.\n\n$d = $word.ActiveDocument\n\n# The target is the last (trailing empty) paragraph in the document body,\n# immediately before the section properties.\n$r = $d.Paragraphs.Last.Range\n$r.InsertAfter(\"<!-- https://mvnrepository.com/artifact/com.github.librepdf/openpdf -->\")\n\n# Each subsequent line becomes its own new paragraph inserted after the\n# previous one, matching the diff's five inserted <w:p> blocks.\n$lines = @(\n    \"<dependency>\",\n    \"    <groupId>com.github.librepdf</groupId>\",\n    \"    <artifactId>openpdf</artifactId>\",\n    \"    <version>1.3.25</version>\",\n    \"</dependency>\"\n)\n\nforeach ($line in $lines) {\n    $r = $d.Paragraphs.Last.Range\n    $r.InsertParagraphAfter()\n    $r = $d.Paragraphs.Last.Range\n    $r.InsertAfter($line)\n}\n"}
